$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.20950566666667
$ws.Range("H2").Value = 87.628517
$ws.Range("I2").Value = 0.01829497698069002
$ws.Range("J2").Value = 0.01840828041918582
$ws.Range("M2").Value = 3.031236
$ws.Range("N2").Value = 9.093708000000001
$ws.Range("O2").Value = 0.6569357730026921
$ws.Range("P2").Value = 0.6780946256479073
$ws.Range("Q2").Value = 88.54090511900401
$ws.Range("R2").Value = 796.8681460710361
$ws.Range("S2").Value = 0.01201862484487606
$ws.Range("T2").Value = 0.01248255601966951
$ws.Range("G3").Value = 29.20950566666667
$ws.Range("H3").Value = 87.628517
$ws.Range("I3").Value = 0.01829497698069002
$ws.Range("J3").Value = 0.01840828041918582
$ws.Range("O3").Value = 0.2314409052885859
$ws.Range("P3").Value = 0.2388952474211406
$ws.Range("Q3").Value = 31.19328871702211
$ws.Range("R3").Value = 280.739598453199
$ws.Range("S3").Value = 0.004234206034644738
$ws.Range("T3").Value = 0.004397650705339135
$ws.Range("G4").Value = 29.20950566666667
$ws.Range("H4").Value = 87.628517
$ws.Range("I4").Value = 0.01829497698069002
$ws.Range("J4").Value = 0.01840828041918582
$ws.Range("M4").Value = 0.04253966666666667
$ws.Range("N4").Value = 0.127619
$ws.Range("O4").Value = 0.009219285072143351
$ws.Range("P4").Value = 0.009516223528461688
$ws.Range("Q4").Value = 1.242562634558111
$ws.Range("R4").Value = 11.183063711023
$ws.Range("S4").Value = 0.0001686666081732817
$ws.Range("T4").Value = 0.0001751773112435767
$ws.Range("G5").Value = 29.20950566666667
$ws.Range("H5").Value = 87.628517
$ws.Range("I5").Value = 0.01829497698069002
$ws.Range("J5").Value = 0.01840828041918582
$ws.Range("M5").Value = 0.4319364999999999
$ws.Range("N5").Value = 0.8638729999999999
$ws.Range("O5").Value = 0.09361017700884301
$ws.Range("P5").Value = 0.06441680759293508
$ws.Range("Q5").Value = 12.61665164439017
$ws.Range("R5").Value = 75.69990986634099
$ws.Range("S5").Value = 0.001712596033535101
$ws.Range("T5").Value = 0.001185802657879487
$ws.Range("G6").Value = 29.20950566666667
$ws.Range("H6").Value = 87.628517
$ws.Range("I6").Value = 0.01829497698069002
$ws.Range("J6").Value = 0.01840828041918582
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.04057666666666667
$ws.Range("N6").Value = 0.12173
$ws.Range("O6").Value = 0.008793859627735762
$ws.Range("P6").Value = 0.009077095809555328
$ws.Range("Q6").Value = 1.185224374934445
$ws.Range("R6").Value = 10.66701937441
$ws.Range("S6").Value = 0.0001608834594608451
$ws.Range("T6").Value = 0.000167093725054111
$ws.Range("I7").Value = 0.913374480506715
$ws.Range("J7").Value = 0.9190311407684336
$ws.Range("M7").Value = 3.031236
$ws.Range("N7").Value = 9.093708000000001
$ws.Range("O7").Value = 0.6569357730026921
$ws.Range("P7").Value = 0.6780946256479073
$ws.Range("Q7").Value = 4420.393821868284
$ws.Range("R7").Value = 39783.54439681456
$ws.Range("S7").Value = 0.6000283703926111
$ws.Range("T7").Value = 0.6231900773581402
$ws.Range("I8").Value = 0.913374480506715
$ws.Range("J8").Value = 0.9190311407684336
$ws.Range("O8").Value = 0.2314409052885859
$ws.Range("P8").Value = 0.2388952474211406
$ws.Range("S8").Value = 0.2113922166359659
$ws.Range("T8").Value = 0.2195521717616081
$ws.Range("I9").Value = 0.913374480506715
$ws.Range("J9").Value = 0.9190311407684336
$ws.Range("M9").Value = 0.04253966666666667
$ws.Range("N9").Value = 0.127619
$ws.Range("O9").Value = 0.009219285072143351
$ws.Range("P9").Value = 0.009516223528461688
$ws.Range("Q9").Value = 62.03478703659811
$ws.Range("R9").Value = 558.3130833293831
$ws.Range("S9").Value = 0.008420659713412246
$ws.Range("T9").Value = 0.008745705765169554
$ws.Range("I10").Value = 0.913374480506715
$ws.Range("J10").Value = 0.9190311407684336
$ws.Range("M10").Value = 0.4319364999999999
$ws.Range("N10").Value = 0.8638729999999999
$ws.Range("O10").Value = 0.09361017700884301
$ws.Range("P10").Value = 0.06441680759293508
$ws.Range("Q10").Value = 629.8847849654101
$ws.Range("R10").Value = 3779.30870979246
$ws.Range("S10").Value = 0.08550114679559362
$ws.Range("T10").Value = 0.05920105216679582
$ws.Range("I11").Value = 0.913374480506715
$ws.Range("J11").Value = 0.9190311407684336
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.04057666666666667
$ws.Range("N11").Value = 0.12173
$ws.Range("O11").Value = 0.008793859627735762
$ws.Range("P11").Value = 0.009077095809555328
$ws.Range("Q11").Value = 59.17218146173445
$ws.Range("R11").Value = 532.54963315561
$ws.Range("S11").Value = 0.008032086969132126
$ws.Range("T11").Value = 0.008342133716720001
$ws.Range("G12").Value = 57.98602933333333
$ws.Range("H12").Value = 173.958088
$ws.Range("I12").Value = 0.03631876156896331
$ws.Range("J12").Value = 0.03654368891224535
$ws.Range("M12").Value = 3.031236
$ws.Range("N12").Value = 9.093708000000001
$ws.Range("O12").Value = 0.6569357730026921
$ws.Range("P12").Value = 0.6780946256479073
$ws.Range("Q12").Value = 175.769339612256
$ws.Range("R12").Value = 1581.924056510304
$ws.Range("S12").Value = 0.02385909370580738
$ws.Range("T12").Value = 0.02478007905274259
$ws.Range("G13").Value = 57.98602933333333
$ws.Range("H13").Value = 173.958088
$ws.Range("I13").Value = 0.03631876156896331
$ws.Range("J13").Value = 0.03654368891224535
$ws.Range("O13").Value = 0.2314409052885859
$ws.Range("P13").Value = 0.2388952474211406
$ws.Range("Q13").Value = 61.92418917285955
$ws.Range("R13").Value = 557.3177025557359
$ws.Range("S13").Value = 0.008405647056481171
$ws.Range("T13").Value = 0.008730113604372045
$ws.Range("G14").Value = 57.98602933333333
$ws.Range("H14").Value = 173.958088
$ws.Range("I14").Value = 0.03631876156896331
$ws.Range("J14").Value = 0.03654368891224535
$ws.Range("M14").Value = 0.04253966666666667
$ws.Range("N14").Value = 0.127619
$ws.Range("O14").Value = 0.009219285072143351
$ws.Range("P14").Value = 0.009516223528461688
$ws.Range("Q14").Value = 2.466706359163556
$ws.Range("R14").Value = 22.200357232472
$ws.Range("S14").Value = 0.0003348330163714771
$ws.Range("T14").Value = 0.0003477579122434937
$ws.Range("G15").Value = 57.98602933333333
$ws.Range("H15").Value = 173.958088
$ws.Range("I15").Value = 0.03631876156896331
$ws.Range("J15").Value = 0.03654368891224535
$ws.Range("M15").Value = 0.4319364999999999
$ws.Range("N15").Value = 0.8638729999999999
$ws.Range("O15").Value = 0.09361017700884301
$ws.Range("P15").Value = 0.06441680759293508
$ws.Range("Q15").Value = 25.04628255913733
$ws.Range("R15").Value = 150.2776953548239
$ws.Range("S15").Value = 0.003399805699212621
$ws.Range("T15").Value = 0.002354027777396184
$ws.Range("G16").Value = 57.98602933333333
$ws.Range("H16").Value = 173.958088
$ws.Range("I16").Value = 0.03631876156896331
$ws.Range("J16").Value = 0.03654368891224535
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.04057666666666667
$ws.Range("N16").Value = 0.12173
$ws.Range("O16").Value = 0.008793859627735762
$ws.Range("P16").Value = 0.009077095809555328
$ws.Range("Q16").Value = 2.352879783582222
$ws.Range("R16").Value = 21.17591805224
$ws.Range("S16").Value = 0.0003193820910906676
$ws.Range("T16").Value = 0.0003317105654910357
$ws.Range("G17").Value = 29.481085
$ws.Range("H17").Value = 58.96217
$ws.Range("I17").Value = 0.01846507700595112
$ws.Range("J17").Value = 0.01238628926567028
$ws.Range("M17").Value = 3.031236
$ws.Range("N17").Value = 9.093708000000001
$ws.Range("O17").Value = 0.6569357730026921
$ws.Range("P17").Value = 0.6780946256479073
$ws.Range("Q17").Value = 89.36412617106001
$ws.Range("R17").Value = 536.1847570263601
$ws.Range("S17").Value = 0.01213036963645874
$ws.Range("T17").Value = 0.00839907618277138
$ws.Range("G18").Value = 29.481085
$ws.Range("H18").Value = 58.96217
$ws.Range("I18").Value = 0.01846507700595112
$ws.Range("J18").Value = 0.01238628926567028
$ws.Range("O18").Value = 0.2314409052885859
$ws.Range("P18").Value = 0.2388952474211406
$ws.Range("Q18").Value = 31.48331254183167
$ws.Range("R18").Value = 188.89987525099
$ws.Range("S18").Value = 0.004273574138480779
$ws.Range("T18").Value = 0.002959025638752119
$ws.Range("G19").Value = 29.481085
$ws.Range("H19").Value = 58.96217
$ws.Range("I19").Value = 0.01846507700595112
$ws.Range("J19").Value = 0.01238628926567028
$ws.Range("M19").Value = 0.04253966666666667
$ws.Range("N19").Value = 0.127619
$ws.Range("O19").Value = 0.009219285072143351
$ws.Range("P19").Value = 0.009516223528461688
$ws.Range("Q19").Value = 1.254115528871667
$ws.Range("R19").Value = 7.524693173230001
$ws.Range("S19").Value = 0.0001702348087969427
$ws.Range("T19").Value = 0.0001178706973403039
$ws.Range("G20").Value = 29.481085
$ws.Range("H20").Value = 58.96217
$ws.Range("I20").Value = 0.01846507700595112
$ws.Range("J20").Value = 0.01238628926567028
$ws.Range("M20").Value = 0.4319364999999999
$ws.Range("N20").Value = 0.8638729999999999
$ws.Range("O20").Value = 0.09361017700884301
$ws.Range("P20").Value = 0.06441680759293508
$ws.Range("Q20").Value = 12.7339566711025
$ws.Range("R20").Value = 50.93582668441
$ws.Range("S20").Value = 0.001728519127009002
$ws.Range("T20").Value = 0.0007978852124171194
$ws.Range("G21").Value = 29.481085
$ws.Range("H21").Value = 58.96217
$ws.Range("I21").Value = 0.01846507700595112
$ws.Range("J21").Value = 0.01238628926567028
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.04057666666666667
$ws.Range("N21").Value = 0.12173
$ws.Range("O21").Value = 0.008793859627735762
$ws.Range("P21").Value = 0.009077095809555328
$ws.Range("Q21").Value = 1.196244159016667
$ws.Range("R21").Value = 7.1774649541
$ws.Range("S21").Value = 0.0001623792952056655
$ws.Range("T21").Value = 0.0001124315343893558
$ws.Range("G22").Value = 21.628479
$ws.Range("H22").Value = 64.885437
$ws.Range("I22").Value = 0.01354670393768061
$ws.Range("J22").Value = 0.01363060063446486
$ws.Range("M22").Value = 3.031236
$ws.Range("N22").Value = 9.093708000000001
$ws.Range("O22").Value = 0.6569357730026921
$ws.Range("P22").Value = 0.6780946256479073
$ws.Range("Q22").Value = 65.56102417004401
$ws.Range("R22").Value = 590.049217530396
$ws.Range("S22").Value = 0.008899314422938825
$ws.Range("T22").Value = 0.009242837034583579
$ws.Range("G23").Value = 21.628479
$ws.Range("H23").Value = 64.885437
$ws.Range("I23").Value = 0.01354670393768061
$ws.Range("J23").Value = 0.01363060063446486
$ws.Range("O23").Value = 0.2314409052885859
$ws.Range("P23").Value = 0.2388952474211406
$ws.Range("Q23").Value = 23.097391570271
$ws.Range("R23").Value = 207.876524132439
$ws.Range("S23").Value = 0.003135261423013251
$ws.Range("T23").Value = 0.00325628571106924
$ws.Range("G24").Value = 21.628479
$ws.Range("H24").Value = 64.885437
$ws.Range("I24").Value = 0.01354670393768061
$ws.Range("J24").Value = 0.01363060063446486
$ws.Range("M24").Value = 0.04253966666666667
$ws.Range("N24").Value = 0.127619
$ws.Range("O24").Value = 0.009219285072143351
$ws.Range("P24").Value = 0.009516223528461688
$ws.Range("Q24").Value = 0.9200682871670001
$ws.Range("R24").Value = 8.280614584503001
$ws.Range("S24").Value = 0.0001248909253894044
$ws.Range("T24").Value = 0.0001297118424647594
$ws.Range("G25").Value = 21.628479
$ws.Range("H25").Value = 64.885437
$ws.Range("I25").Value = 0.01354670393768061
$ws.Range("J25").Value = 0.01363060063446486
$ws.Range("M25").Value = 0.4319364999999999
$ws.Range("N25").Value = 0.8638729999999999
$ws.Range("O25").Value = 0.09361017700884301
$ws.Range("P25").Value = 0.06441680759293508
$ws.Range("Q25").Value = 9.342129519583498
$ws.Range("R25").Value = 56.05277711750099
$ws.Range("S25").Value = 0.001268109353492673
$ws.Range("T25").Value = 0.0008780397784464619
$ws.Range("G26").Value = 21.628479
$ws.Range("H26").Value = 64.885437
$ws.Range("I26").Value = 0.01354670393768061
$ws.Range("J26").Value = 0.01363060063446486
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.04057666666666667
$ws.Range("N26").Value = 0.12173
$ws.Range("O26").Value = 0.008793859627735762
$ws.Range("P26").Value = 0.009077095809555328
$ws.Range("Q26").Value = 0.8776115828900001
$ws.Range("R26").Value = 7.89850424601
$ws.Range("S26").Value = 0.0001191278128464586
$ws.Range("T26").Value = 0.0001237262679008232
